# The worksheet is protected; unprotect first so cell values can be edited,
# then restore protection afterwards.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A10).
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-07 for illustrative purposes only and are subject to change."

# Update weight (col D) and percent-change (col E) figures for rows 2-7.
$ws.Range("D2").Value = 0.4994294177633177
$ws.Range("E2").Value = 0

$ws.Range("D3").Value = 0.323130432671296
$ws.Range("E3").Value = 0

$ws.Range("D4").Value = 0.08989705239644538
$ws.Range("E4").Value = 0

$ws.Range("D5").Value = 0.05876450735686883
$ws.Range("E5").Value = 0

$ws.Range("D6").Value = 0.02877858981207206
$ws.Range("E6").Value = 0

$ws.Range("E7").Value = 0

$ws.Protect()
